# Update the dSF column (F) values for several rows to reflect
# re-pulled data / recalculated means.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F2"  = -2
    "F5"  = -1
    "F9"  = 0
    "F10" = 1
    "F11" = 0
    "F14" = -9
    "F15" = -2
    "F20" = -2
    "F21" = -9
    "F22" = -9
    "F23" = 1
    "F27" = 4
    "F31" = 1
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
